$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 08:43"

# Swap Islas Malvinas / Montserrat rows (214/215)
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Row 56 - Kirguistan
$ws.Range("B56").Value = 43245
$ws.Range("C56").Value = 119
$ws.Range("D56").Value = 36925
$ws.Range("E56").Value = 5263

# Row 62 - Uzbekistan
$ws.Range("B62").Value = 39506
$ws.Range("C62").Value = 158
$ws.Range("E62").Value = 3673
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 282

# Row 63 - Afganistan
$ws.Range("B63").Value = 38070
$ws.Range("C63").Value = 16
$ws.Range("D63").Value = 28440
$ws.Range("E63").Value = 8233
$ws.Range("G63").Value = 8
$ws.Range("H63").Value = 1397

# Row 73 - El Salvador
$ws.Range("D73").Value = 12506
$ws.Range("E73").Value = 11627
$ws.Range("G73").Value = 9
$ws.Range("H73").Value = 678

# Row 151 - Georgia
$ws.Range("B151").Value = 1429
$ws.Range("C151").Value = 8
$ws.Range("D151").Value = 1150
$ws.Range("E151").Value = 260
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 19

# Row 161 - Vietnam
$ws.Range("D161").Value = 590
$ws.Range("E161").Value = 405
